# Updated cryptos list on Mon Mar 25 09:21:43 UTC 2024 with GitHub Actions
#
# Refreshes the Price (column D) and Volume(1h) (column E) figures for
# each tracked coin, and reflects two ranking swaps that occurred in
# this run: dogwifhat overtook Fetch.AI (rows 45/46) and LidoDAOToken
# overtook Stellar (rows 48/49).
#
# Note: several Price values are plain decimals (e.g. "579.13") that
# Excel would otherwise auto-convert to numbers on assignment, losing
# the original text formatting (e.g. "18.70" -> 18.7). To preserve the
# text representation exactly as authored, those cells are written
# with a leading apostrophe (forces text entry) and then restored to
# the workbook's default "Normal" style so no extra numeric format
# is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.833.36'
$ws.Range("E2").Value = '  +2.90%  '
$ws.Range("D3").Value = '3.443.95'
$ws.Range("E3").Value = '  +2.29%  '
$ws.Range("E4").Value = '  -0.03%  '
$r = $ws.Range("D5")
$r.Value = "'579.13"
$r.Style = "Normal"
$ws.Range("E5").Value = '  +4.45%  '
$r = $ws.Range("D6")
$r.Value = "'187.33"
$r.Style = "Normal"
$ws.Range("E6").Value = '  +7.59%  '
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("D8").Value = '3.436.58'
$ws.Range("E8").Value = '  +2.27%  '
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("E10").Value = '  -2.04%  '
$r = $ws.Range("D11")
$r.Value = "'0.644"
$r.Style = "Normal"
$ws.Range("E11").Value = '  +1.26%  '
$r = $ws.Range("D12")
$r.Value = "'56.55"
$r.Style = "Normal"
$ws.Range("E12").Value = '  +5.62%  '
$ws.Range("E13").Value = '  -1.61%  '
$r = $ws.Range("D14")
$r.Value = "'9.39"
$r.Style = "Normal"
$ws.Range("E14").Value = '  +2.67%  '
$ws.Range("D15").Value = '3.994.82'
$r = $ws.Range("D16")
$r.Value = "'18.70"
$r.Style = "Normal"
$ws.Range("E16").Value = '  +2.58%  '
$ws.Range("D17").Value = '3.449.60'
$ws.Range("E17").Value = '  +2.70%  '
$ws.Range("D18").Value = '66.897.97'
$ws.Range("E18").Value = '  +2.93%  '
$r = $ws.Range("D19")
$r.Value = "'12.07"
$r.Style = "Normal"
$ws.Range("E19").Value = '  +1.92%  '
$ws.Range("E20").Value = '  -1.98%  '
$ws.Range("E21").Value = '  +2.42%  '
$r = $ws.Range("D22")
$r.Value = "'482.30"
$r.Style = "Normal"
$ws.Range("E22").Value = '  +6.80%  '
$r = $ws.Range("D23")
$r.Value = "'5.32"
$r.Style = "Normal"
$ws.Range("E23").Value = '  +7.80%  '
$r = $ws.Range("D24")
$r.Value = "'16.86"
$r.Style = "Normal"
$ws.Range("E24").Value = '  +22.99%  '
$r = $ws.Range("D25")
$r.Value = "'4.33"
$r.Style = "Normal"
$ws.Range("E25").Value = '  +6.39%  '
$r = $ws.Range("D26")
$r.Value = "'89.34"
$r.Style = "Normal"
$ws.Range("E26").Value = '  +2.79%  '
$ws.Range("E27").Value = '  +3.04%  '
$r = $ws.Range("D28")
$r.Value = "'10.97"
$r.Style = "Normal"
$ws.Range("E28").Value = '  +2.02%  '
$r = $ws.Range("D29")
$r.Value = "'9.04"
$r.Style = "Normal"
$ws.Range("E29").Value = '  +4.68%  '
$r = $ws.Range("D30")
$r.Value = "'31.22"
$r.Style = "Normal"
$ws.Range("E30").Value = '  +0.70%  '
$r = $ws.Range("D31")
$r.Value = "'7.29"
$r.Style = "Normal"
$ws.Range("E31").Value = '  +11.24%  '
$ws.Range("E32").Value = '  +2.68%  '
$r = $ws.Range("D33")
$r.Value = "'597.08"
$r.Style = "Normal"
$ws.Range("E33").Value = '  +3.63%  '
$r = $ws.Range("D34")
$r.Value = "'63.63"
$r.Style = "Normal"
$ws.Range("E34").Value = '  +1.16%  '
$r = $ws.Range("D35")
$r.Value = "'0.111"
$r.Style = "Normal"
$ws.Range("E35").Value = '  +3.77%  '
$r = $ws.Range("D36")
$r.Value = "'0.149"
$r.Style = "Normal"
$ws.Range("E36").Value = '  +6.34%  '
$ws.Range("E37").Value = '  -0.07%  '
$r = $ws.Range("D38")
$r.Value = "'36.64"
$r.Style = "Normal"
$ws.Range("E38").Value = '  +2.88%  '
$r = $ws.Range("D39")
$r.Value = "'0.387"
$r.Style = "Normal"
$ws.Range("E39").Value = '  +4.32%  '
$r = $ws.Range("D40")
$r.Value = "'3.53"
$r.Style = "Normal"
$ws.Range("E40").Value = '  -2.62%  '
$ws.Range("D41").Value = '3.249.16'
$ws.Range("E41").Value = '  +5.22%  '
$ws.Range("E42").Value = '  +1.72%  '
$r = $ws.Range("D43")
$r.Value = "'2.89"
$r.Style = "Normal"
$ws.Range("E43").Value = '  +4.64%  '
$r = $ws.Range("D44")
$r.Value = "'0.0429"
$r.Style = "Normal"
$ws.Range("E44").Value = '  +2.89%  '
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$r = $ws.Range("D45")
$r.Value = "'2.84"
$r.Style = "Normal"
$ws.Range("E45").Value = '  +25.44%  '
$ws.Range("B46").Value = 'Fetch.AI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$r = $ws.Range("D46")
$r.Value = "'2.53"
$r.Style = "Normal"
$ws.Range("E46").Value = '  +3.40%  '
$ws.Range("E47").Value = '  +1.80%  '
$ws.Range("B48").Value = 'LidoDAOToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$r = $ws.Range("D48")
$r.Value = "'3.36"
$r.Style = "Normal"
$ws.Range("E48").Value = '  +15.72%  '
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$r = $ws.Range("D49")
$r.Value = "'0.134"
$r.Style = "Normal"
$ws.Range("E49").Value = '  +0.45%  '
$ws.Range("E50").Value = '  +0.04%  '
$r = $ws.Range("D51")
$r.Value = "'8.69"
$r.Style = "Normal"
$ws.Range("E51").Value = '  +5.18%  '
